# Apply updated team-specific transition matrix values
# (changes to team matrices from games pulled march 7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1976744186046512
$ws.Range("C2").Value = 0.5348837209302325
$ws.Range("J2").Value = 0.01453488372093023
$ws.Range("P2").Value = 0.1308139534883721
$ws.Range("S2").Value = 0.1220930232558139
$ws.Range("B3").Value = 0.02030456852791878
$ws.Range("C3").Value = 0.02538071065989848
$ws.Range("J3").Value = 0.03553299492385787
$ws.Range("P3").Value = 0.6852791878172588
$ws.Range("S3").Value = 0.233502538071066
$ws.Range("J4").Value = 0.03773584905660377
$ws.Range("O4").Value = 0.01886792452830189
$ws.Range("P4").Value = 0.660377358490566
$ws.Range("S4").Value = 0.2830188679245283
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.0546875
$ws.Range("D6").Value = 0.0078125
$ws.Range("F6").Value = 0.1328125
$ws.Range("J6").Value = 0.22265625
$ws.Range("O6").Value = 0.03125
$ws.Range("Q6").Value = 0.16015625
$ws.Range("R6").Value = 0.05859375
$ws.Range("S6").Value = 0.33203125
$ws.Range("B7").Value = 0.0958904109589041
$ws.Range("D7").Value = 0.0045662100456621
$ws.Range("F7").Value = 0.0502283105022831
$ws.Range("J7").Value = 0.1187214611872146
$ws.Range("O7").Value = 0.0273972602739726
$ws.Range("Q7").Value = 0.1872146118721461
$ws.Range("R7").Value = 0.136986301369863
$ws.Range("S7").Value = 0.3789954337899543
$ws.Range("B8").Value = 0.1026200873362445
$ws.Range("D8").Value = 0.0240174672489083
$ws.Range("E8").Value = 0.002183406113537118
$ws.Range("F8").Value = 0.05458515283842795
$ws.Range("J8").Value = 0.1091703056768559
$ws.Range("O8").Value = 0.03056768558951965
$ws.Range("Q8").Value = 0.1812227074235808
$ws.Range("R8").Value = 0.1135371179039301
$ws.Range("S8").Value = 0.3820960698689956
$ws.Range("B9").Value = 0.0846774193548387
$ws.Range("D9").Value = 0.01209677419354839
$ws.Range("F9").Value = 0.05241935483870968
$ws.Range("J9").Value = 0.1532258064516129
$ws.Range("O9").Value = 0.01209677419354839
$ws.Range("Q9").Value = 0.1774193548387097
$ws.Range("R9").Value = 0.08870967741935484
$ws.Range("S9").Value = 0.4193548387096774
$ws.Range("B10").Value = 0.1207764198418404
$ws.Range("D10").Value = 0.02659956865564342
$ws.Range("E10").Value = 0.001437814521926672
$ws.Range("F10").Value = 0.06470165348670022
$ws.Range("J10").Value = 0.1157440690150971
$ws.Range("O10").Value = 0.01581595974119339
$ws.Range("Q10").Value = 0.2135154565061107
$ws.Range("R10").Value = 0.07548526240115025
$ws.Range("S10").Value = 0.3659237958303379
$ws.Range("G11").Value = 0.1564417177914111
$ws.Range("J11").Value = 0.06748466257668712
$ws.Range("K11").Value = 0.1993865030674846
$ws.Range("L11").Value = 0.5552147239263804
$ws.Range("S11").Value = 0.02147239263803681
$ws.Range("G12").Value = 0.7409326424870466
$ws.Range("J12").Value = 0.1450777202072539
$ws.Range("K12").Value = 0.02590673575129534
$ws.Range("L12").Value = 0.05699481865284974
$ws.Range("S12").Value = 0.0310880829015544
$ws.Range("G13").Value = 0.68
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.02
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02325581395348837
$ws.Range("H15").Value = 0.1356589147286822
$ws.Range("I15").Value = 0.06589147286821706
$ws.Range("J15").Value = 0.3023255813953488
$ws.Range("K15").Value = 0.09689922480620156
$ws.Range("M15").Value = 0.003875968992248062
$ws.Range("O15").Value = 0.1124031007751938
$ws.Range("S15").Value = 0.2596899224806202
$ws.Range("F16").Value = 0.01428571428571429
$ws.Range("H16").Value = 0.1476190476190476
$ws.Range("I16").Value = 0.08571428571428572
$ws.Range("J16").Value = 0.4142857142857143
$ws.Range("K16").Value = 0.09047619047619047
$ws.Range("M16").Value = 0.01904761904761905
$ws.Range("O16").Value = 0.06190476190476191
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("F17").Value = 0.02208835341365462
$ws.Range("H17").Value = 0.2048192771084337
$ws.Range("I17").Value = 0.09839357429718876
$ws.Range("J17").Value = 0.3975903614457831
$ws.Range("K17").Value = 0.08835341365461848
$ws.Range("M17").Value = 0.02208835341365462
$ws.Range("O17").Value = 0.06224899598393574
$ws.Range("S17").Value = 0.1044176706827309
$ws.Range("F18").Value = 0.01339285714285714
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1294642857142857
$ws.Range("J18").Value = 0.40625
$ws.Range("K18").Value = 0.07589285714285714
$ws.Range("M18").Value = 0.01785714285714286
$ws.Range("O18").Value = 0.08482142857142858
$ws.Range("S18").Value = 0.1294642857142857
$ws.Range("F19").Value = 0.0189873417721519
$ws.Range("H19").Value = 0.1849507735583685
$ws.Range("I19").Value = 0.0949367088607595
$ws.Range("J19").Value = 0.3783403656821379
$ws.Range("K19").Value = 0.1061884669479606
$ws.Range("M19").Value = 0.02320675105485232
$ws.Range("N19").Value = 0.002109704641350211
$ws.Range("O19").Value = 0.06118143459915612
$ws.Range("S19").Value = 0.130098452883263
